# Word COM-interop script implementing the commit:
# "k means clustering applied to vectorized and LSA vectorized feature arrays,
#  the vectorized array got much better result than LSA reduced array"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: add a new paragraph right after the paragraph that ends with
# "...which process can take minutes." describing the K-means++ clustering
# result.
# ---------------------------------------------------------------------------
$i = 1
$targetIndex = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*which process can take minutes.*") {
        $targetIndex = $i
    }
    $i = $i + 1
}

$target = $d.Paragraphs($targetIndex)
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($targetIndex + 1)
$newPara.Range.Text = "Using K means++ to cluster, the vectorized 5000- feature array got much better clustering resulst than using LSA reduced 2 feature array. Maybe increasing features will get better result?"

# ---------------------------------------------------------------------------
# Change 2: fix the typo "labes" -> "labels" and add the "of original"
# qualifier in the "Add label text files..." to-do item.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("labes and km.labels", $true, $false, $false, $false, $false,
                         $true, 1, $false, "labels of original and km.labels", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3 (side effect of continued editing of item 3): the "_GoBack" edit
# bookmark relocates from item 2 to the middle of item 3, between
# "...methods-o" and "pen to fix and try".
# ---------------------------------------------------------------------------
try {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldBookmark.Delete()
} catch {
}
$bmFind = $d.Content
$bmFound = $bmFind.Find.Execute("these two feature extraction methods" + [char]0x2014 + "o")
if ($bmFound) {
    $bmSpot = $d.Range($bmFind.End, $bmFind.End)
    $d.Bookmarks.Add("_GoBack", $bmSpot) | Out-Null
}

# ---------------------------------------------------------------------------
# Change 4: append a new to-do item ("4. Cross - validation needed? ...")
# right after the "...TfidfVectorizer" paragraph (item 3), before the
# trailing empty paragraph at the end of the document.
# ---------------------------------------------------------------------------
$i = 1
$idx3 = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*if these 2 methods will work better than TfidfVectorizer*") {
        $idx3 = $i
    }
    $i = $i + 1
}

$p3 = $d.Paragraphs($idx3)
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs($idx3 + 1)
$p4.Range.Text = "4. Cross " + [char]0x2013 + " validation needed? How to detect suspicious outlier( doubted not Shakespear authored works?)"
